# Update the cryptos price/volume list (GitHub Actions refresh).
# For Price (column D) cells whose new text would otherwise be auto-detected
# as a number by Excel (losing trailing zeros / exact formatting), force the
# cell to Text format first so the literal string is preserved, matching the
# source workbook's inline-string cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.683.46'
$ws.Range('E2').Value = '  -2.10%  '
$ws.Range('D3').Value = '1.797.69'
$ws.Range('E3').Value = '  -1.85%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '231.69'
$ws.Range('E5').Value = '  -1.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5881'
$ws.Range('E6').Value = '  -2.52%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  -0.65%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06767'
$ws.Range('E9').Value = '  -4.06%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.23'
$ws.Range('E10').Value = '  -1.26%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07527'
$ws.Range('D12').Value = '1.794.79'
$ws.Range('E12').Value = '  -1.98%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.790'
$ws.Range('E13').Value = '  -0.15%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6178'
$ws.Range('E14').Value = '  -1.71%  '
$ws.Range('D15').Value = '2.041.18'
$ws.Range('E15').Value = '  -1.86%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000009101'
$ws.Range('E16').Value = '  -8.30%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '75.31'
$ws.Range('E17').Value = '  -4.88%  '
$ws.Range('D18').Value = '28.653.34'
$ws.Range('E18').Value = '  -2.15%  '
$ws.Range('E19').Value = '  -6.38%  '
$ws.Range('E20').Value = '  -0.10%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '210.52'
$ws.Range('E21').Value = '  -6.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.50'
$ws.Range('E22').Value = '  -1.77%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.827'
$ws.Range('E24').Value = '  -0.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.61'
$ws.Range('E25').Value = '  -1.45%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.073'
$ws.Range('E26').Value = '  +1.10%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1262'
$ws.Range('E27').Value = '  -2.86%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.46'
$ws.Range('E28').Value = '  -0.69%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.425'
$ws.Range('E29').Value = '  -3.57%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.06110'
$ws.Range('E30').Value = '  -3.60%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.423'
$ws.Range('E31').Value = '  -1.78%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.807'
$ws.Range('E32').Value = '  +0.05%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.793'
$ws.Range('E33').Value = '  -1.41%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.736'
$ws.Range('E34').Value = '  +0.08%  '
$ws.Range('E35').Value = '  -5.39%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6427'
$ws.Range('E36').Value = '  -0.54%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.500'
$ws.Range('E37').Value = '  -1.99%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.713'
$ws.Range('E38').Value = '  -1.24%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.459'
$ws.Range('E39').Value = '  -0.80%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01701'
$ws.Range('E40').Value = '  -2.09%  '
$ws.Range('D41').Value = '1.142.50'
$ws.Range('E41').Value = '  -6.34%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8839'
$ws.Range('E42').Value = '  -1.68%  '
$ws.Range('E43').Value = '  +0.16%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.22'
$ws.Range('E44').Value = '  -0.32%  '
$ws.Range('D45').Value = '1.950.16'
$ws.Range('E45').Value = '  -2.17%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '60.07'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000112'
$ws.Range('E47').Value = '  -4.12%  '
$ws.Range('E48').Value = '  +0.67%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05489'
$ws.Range('E49').Value = '  -0.20%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.358'
$ws.Range('E50').Value = '  -2.28%  '
$ws.Range('E51').Value = '  -1.82%  '
